$wb = $excel.ActiveWorkbook
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsTasas = $wb.Worksheets.Item("tasas")

# Update the "Conversión del día" note text with the new rates (Hoja1!A1)
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.87 = 6779.69 pesos`n✅ 6779.69 pesos = 1.86 = 966.67 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# Update the numeric rate cells on the "tasas" sheet
$wsTasas.Range("N10").Value = 534.98
$wsTasas.Range("O10").Value = 3627
$wsTasas.Range("N12").Value = 3639.98
$wsTasas.Range("O12").Value = 519
